# Refresh the cryptocurrency price/volume snapshot (and restore the correct
# OKB/Dogecoin and Filecoin/InternetComputer(DFINITY) row ordering) to match
# the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''29.023.06'
$ws.Range("E2").Value = '  -0.46%  '

# Row 3
$ws.Range("D3").Value = '''1.828.89'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("D4").Value = '''0.9992'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '''241.02'
$ws.Range("E5").Value = '  -0.21%  '

# Row 6
$ws.Range("D6").Value = '''0.6225'
$ws.Range("E6").Value = '  -5.46%  '

# Row 7
$ws.Range("D7").Value = '''1.001'

# Row 8
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '''44.60'
$ws.Range("E8").Value = '  +6.72%  '

# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.07517'
$ws.Range("E9").Value = '  +1.66%  '

# Row 10
$ws.Range("D10").Value = '''0.2906'
$ws.Range("E10").Value = '  -0.62%  '

# Row 11
$ws.Range("D11").Value = '''22.73'

# Row 12
$ws.Range("D12").Value = '''0.07637'
$ws.Range("E12").Value = '  -1.47%  '

# Row 13
$ws.Range("D13").Value = '''1.830.43'
$ws.Range("E13").Value = '  -0.15%  '

# Row 14
$ws.Range("D14").Value = '''4.951'
$ws.Range("E14").Value = '  -0.81%  '

# Row 15
$ws.Range("D15").Value = '''0.6633'
$ws.Range("E15").Value = '  -0.23%  '

# Row 16
$ws.Range("D16").Value = '''82.17'
$ws.Range("E16").Value = '  -0.71%  '

# Row 17
$ws.Range("D17").Value = '''0.000009073'
$ws.Range("E17").Value = '  +7.67%  '

# Row 18
$ws.Range("D18").Value = '''5.995'
$ws.Range("E18").Value = '  -1.82%  '

# Row 19
$ws.Range("D19").Value = '''28.797.96'
$ws.Range("E19").Value = '  -1.27%  '

# Row 20
$ws.Range("D20").Value = '''224.26'
$ws.Range("E20").Value = '  -1.17%  '

# Row 21
$ws.Range("E21").Value = '  -0.95%  '

# Row 22
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.01%  '

# Row 23
$ws.Range("D23").Value = '''7.179'
$ws.Range("E23").Value = '  +0.84%  '

# Row 24
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("E25").Value = '  +0.40%  '

# Row 26
$ws.Range("D26").Value = '''8.377'
$ws.Range("E26").Value = '  -2.56%  '

# Row 27
$ws.Range("D27").Value = '''0.1352'
$ws.Range("E27").Value = '  -2.69%  '

# Row 28
$ws.Range("E28").Value = '  -0.53%  '

# Row 29
$ws.Range("D29").Value = '''1.495'
$ws.Range("E29").Value = '  -1.48%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''4.027'
$ws.Range("E30").Value = '  -0.32%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''4.044'
$ws.Range("E31").Value = '  -1.60%  '

# Row 32
$ws.Range("E32").Value = '  +0.96%  '

# Row 33
$ws.Range("D33").Value = '''0.05182'
$ws.Range("E33").Value = '  -1.28%  '

# Row 34
$ws.Range("D34").Value = '''1.833'
$ws.Range("E34").Value = '  -1.57%  '

# Row 35
$ws.Range("E35").Value = '  +0.83%  '

# Row 36
$ws.Range("D36").Value = '''0.7307'
$ws.Range("E36").Value = '  -1.16%  '

# Row 37
$ws.Range("D37").Value = '''2.611'
$ws.Range("E37").Value = '  -1.65%  '

# Row 38
$ws.Range("D38").Value = '''1.284.16'
$ws.Range("E38").Value = '  -1.11%  '

# Row 39
$ws.Range("E39").Value = '  +0.88%  '

# Row 40
$ws.Range("E40").Value = '  -0.62%  '

# Row 41
$ws.Range("D41").Value = '''6.387'
$ws.Range("E41").Value = '  +7.27%  '

# Row 42
$ws.Range("D42").Value = '''0.8955'
$ws.Range("E42").Value = '  -2.62%  '

# Row 43
$ws.Range("E43").Value = '  +0.17%  '

# Row 44
$ws.Range("D44").Value = '''101.40'
$ws.Range("E44").Value = '  -0.84%  '

# Row 45
$ws.Range("D45").Value = '''1.980.17'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("D46").Value = '''0.5117'
$ws.Range("E46").Value = '  -0.49%  '

# Row 47
$ws.Range("D47").Value = '''63.40'
$ws.Range("E47").Value = '  +0.22%  '

# Row 48
$ws.Range("E48").Value = '  -0.13%  '

# Row 49
$ws.Range("D49").Value = '''0.3967'
$ws.Range("E49").Value = '  -0.83%  '

# Row 50
$ws.Range("D50").Value = '''8.842'
$ws.Range("E50").Value = '  +1.32%  '

# Row 51
$ws.Range("D51").Value = '''1.646'
$ws.Range("E51").Value = '  -5.97%  '

# The apostrophe-forced entries above mark D2:D51 with a "quote prefix" style
# (Excel's normal behavior for typed text-that-looks-like-a-number). Clear the
# formatting back off so the cells keep their original unstyled look.
$ws.Range("D2:D51").ClearFormats()
